# Daily attendance processing - 2026-01-30 11:22:27
# Rotate the "Recorded By" (column G) list of recorders so the first
# recorder moves to the end of the comma-separated list.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = "Recorded By"
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*, *") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $rotated = $parts[1..($parts.Length - 1)] + $parts[0]
            $newVal = $rotated -join ", "
            if ($newVal -ne $val) {
                $cell.Value = $newVal
            }
        }
    }
}
